# Update crypto price/volume figures per latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.094.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.49%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5201"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.59%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2622"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.51%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06301"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.38%  "

$ws.Range("E11").Value = "  -1.71%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.584"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.95%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.676.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.875.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.67%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5573"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8116"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.91%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.09%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.063.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.56%  "

$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("E20").Value = "  -2.42%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "194.45"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.18%  "

$ws.Range("E22").Value = "  +2.34%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.919"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.65%  "

$ws.Range("E24").Value = "  -0.15%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1181"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.192"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.31%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.85"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.02%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.512"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.81%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05433"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.80%  "

$ws.Range("E31").Value = "  -0.35%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.434"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.26%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.325"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.98%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.556"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.416"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.52%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.780"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.90%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9413"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5603"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.51%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01570"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.91%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.736"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.28%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.027.80"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.62%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8233"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.67%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.69"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.786.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.72%  "

$ws.Range("E46").Value = "  +7.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.28"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9981"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4321"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.895"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05098"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.07%  "
